$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update SKU names (shared strings content change reflected through displayed cell text)
$ws.Range("B4").Value = "test SKU 301"
$ws.Range("C4").Value = "test SKU 302"
$ws.Range("D4").Value = "test SKU 303"

# B1: 1 -> 3 (NoOfSku)
$ws.Range("B1").Value = 3

# B29/B30 counts
$ws.Range("B29").Value = 6
$ws.Range("B30").Value = 3

# Rows 32-37: SKU text updates
$ws.Range("B32").Value = "test SKU 301"
$ws.Range("H32").Value = "test SKU 301"
$ws.Range("B33").Value = "test SKU 301"
$ws.Range("H33").Value = "test SKU 302"
$ws.Range("B34").Value = "test SKU 302"
$ws.Range("H34").Value = "test SKU 303"
$ws.Range("B35").Value = "test SKU 302"
$ws.Range("B36").Value = "test SKU 303"
$ws.Range("B37").Value = "test SKU 303"

# Update selection
$ws.Range("B1").Select()
